$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Bdnf"
$ws.Range("C2").Value = "Ngfr"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.08140533333333333
$ws.Range("H2").Value = 0.244216
$ws.Range("I2").Value = 0.1131514935296598
$ws.Range("J2").Value = 0.1131514935296598
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.3922183333333333
$ws.Range("N2").Value = 1.176655
$ws.Range("O2").Value = 0.06257714082953221
$ws.Range("P2").Value = 0.06257714082953222
$ws.Range("Q2").Value = 0.03192866416444445
$ws.Range("R2").Value = 0.28735797748
$ws.Range("S2").Value = 0.007080696945677423
$ws.Range("T2").Value = 0.007080696945677424

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Bdnf"
$ws.Range("C3").Value = "Ngfr"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.08140533333333333
$ws.Range("H3").Value = 0.244216
$ws.Range("I3").Value = 0.1131514935296598
$ws.Range("J3").Value = 0.1131514935296598
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 5.787580333333334
$ws.Range("N3").Value = 17.362741
$ws.Range("O3").Value = 0.9233893441524432
$ws.Range("P3").Value = 0.9233893441524432
$ws.Range("Q3").Value = 0.4711399062284444
$ws.Range("R3").Value = 4.240259156056
$ws.Range("S3").Value = 0.104482883400222
$ws.Range("T3").Value = 0.104482883400222

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Bdnf"
$ws.Range("C4").Value = "Ngfr"
$ws.Range("D4").Value = "Neutrophils"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.08140533333333333
$ws.Range("H4").Value = 0.244216
$ws.Range("I4").Value = 0.1131514935296598
$ws.Range("J4").Value = 0.1131514935296598
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.034619
$ws.Range("N4").Value = 0.103857
$ws.Range("O4").Value = 0.005523347213187152
$ws.Range("P4").Value = 0.005523347213187152
$ws.Range("Q4").Value = 0.002818171234666667
$ws.Range("R4").Value = 0.025363541112
$ws.Range("S4").Value = 0.0006249749864550104
$ws.Range("T4").Value = 0.0006249749864550103

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Bdnf"
$ws.Range("C5").Value = "Ngfr"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.08140533333333333
$ws.Range("H5").Value = 0.244216
$ws.Range("I5").Value = 0.1131514935296598
$ws.Range("J5").Value = 0.1131514935296598
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.05333966666666667
$ws.Range("N5").Value = 0.160019
$ws.Range("O5").Value = 0.00851016780483737
$ws.Range("P5").Value = 0.008510167804837371
$ws.Range("Q5").Value = 0.004342133344888889
$ws.Range("R5").Value = 0.039079200104
$ws.Range("S5").Value = 0.0009629381973053746
$ws.Range("T5").Value = 0.0009629381973053746

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Bdnf"
$ws.Range("C6").Value = "Ngfr"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.6380313333333333
$ws.Range("H6").Value = 1.914094
$ws.Range("I6").Value = 0.8868485064703402
$ws.Range("J6").Value = 0.8868485064703401
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.3922183333333333
$ws.Range("N6").Value = 1.176655
$ws.Range("O6").Value = 0.06257714082953221
$ws.Range("P6").Value = 0.06257714082953222
$ws.Range("Q6").Value = 0.2502475861744444
$ws.Range("R6").Value = 2.25222827557
$ws.Range("S6").Value = 0.05549644388385479
$ws.Range("T6").Value = 0.05549644388385479

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Bdnf"
$ws.Range("C7").Value = "Ngfr"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.6380313333333333
$ws.Range("H7").Value = 1.914094
$ws.Range("I7").Value = 0.8868485064703402
$ws.Range("J7").Value = 0.8868485064703401
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 5.787580333333334
$ws.Range("N7").Value = 17.362741
$ws.Range("O7").Value = 0.9233893441524432
$ws.Range("P7").Value = 0.9233893441524432
$ws.Range("Q7").Value = 3.692657596850444
$ws.Range("R7").Value = 33.233918371654
$ws.Range("S7").Value = 0.8189064607522212
$ws.Range("T7").Value = 0.8189064607522212

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Bdnf"
$ws.Range("C8").Value = "Ngfr"
$ws.Range("D8").Value = "Neutrophils"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.6380313333333333
$ws.Range("H8").Value = 1.914094
$ws.Range("I8").Value = 0.8868485064703402
$ws.Range("J8").Value = 0.8868485064703401
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.034619
$ws.Range("N8").Value = 0.103857
$ws.Range("O8").Value = 0.005523347213187152
$ws.Range("P8").Value = 0.005523347213187152
$ws.Range("Q8").Value = 0.02208800672866667
$ws.Range("R8").Value = 0.198792060558
$ws.Range("S8").Value = 0.004898372226732141
$ws.Range("T8").Value = 0.004898372226732141

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Bdnf"
$ws.Range("C9").Value = "Ngfr"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.6380313333333333
$ws.Range("H9").Value = 1.914094
$ws.Range("I9").Value = 0.8868485064703402
$ws.Range("J9").Value = 0.8868485064703401
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.05333966666666667
$ws.Range("N9").Value = 0.160019
$ws.Range("O9").Value = 0.00851016780483737
$ws.Range("P9").Value = 0.008510167804837371
$ws.Range("Q9").Value = 0.03403237864288889
$ws.Range("R9").Value = 0.306291407786
$ws.Range("S9").Value = 0.007547229607531994
$ws.Range("T9").Value = 0.007547229607531995
